$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-7
# from serial date 45212 (2023-10-13) to 45221 (2023-10-22)
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45221
}
